$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.962.65"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.595.35"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.52"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.245"
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("E11").Value = "  +3.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.819.23"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.598.92"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.00"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.512"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.964.52"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.02"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "199.80"
$ws.Range("E20").Value = "  +3.93%  "
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("E24").Value = "  +5.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.71"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("E27").Value = "  -8.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.07"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("E33").Value = "  -2.48%  "
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("E35").Value = "  +2.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.123.31"
$ws.Range("E36").Value = "  +2.10%  "
$ws.Range("E37").Value = "  +7.13%  "
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.31"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.783"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("E41").Value = "  -3.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.780"
$ws.Range("E42").Value = "  -3.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.730.29"
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.54"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0503"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₇0916"
$ws.Range("E51").Value = "  -13.01%  "
